# edit.ps1
# Applies the "Updated cryptos list" data refresh to the worksheet.
# For each changed cell (B/C/D/E columns, rows 2-51) we set the new
# literal value. Columns D/E store plain text (prices/volumes as
# formatted strings, e.g. "29.493.76" or "  +0.19%  "), so for column D
# values that Excel's automatic type inference would otherwise parse as
# a genuine number (e.g. "1.003"), we briefly force a Text number format
# before assigning the value and then restore the cell's style back to
# "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Addr = 'D2'; Value = '29.493.76'; ForceText = 1 },
    @{ Addr = 'E2'; Value = '  +0.19%  '; ForceText = 0 },
    @{ Addr = 'D3'; Value = '1.898.92'; ForceText = 1 },
    @{ Addr = 'E3'; Value = '  -1.00%  '; ForceText = 0 },
    @{ Addr = 'E4'; Value = '  -0.50%  '; ForceText = 0 },
    @{ Addr = 'D5'; Value = '338.51'; ForceText = 1 },
    @{ Addr = 'E5'; Value = '  +4.05%  '; ForceText = 0 },
    @{ Addr = 'D6'; Value = '1.003'; ForceText = 1 },
    @{ Addr = 'D7'; Value = '0.4751'; ForceText = 1 },
    @{ Addr = 'E7'; Value = '  -1.48%  '; ForceText = 0 },
    @{ Addr = 'D8'; Value = '0.3998'; ForceText = 1 },
    @{ Addr = 'E8'; Value = '  -1.72%  '; ForceText = 0 },
    @{ Addr = 'D9'; Value = '0.08035'; ForceText = 1 },
    @{ Addr = 'E9'; Value = '  -2.22%  '; ForceText = 0 },
    @{ Addr = 'D10'; Value = '0.9905'; ForceText = 1 },
    @{ Addr = 'E10'; Value = '  -2.23%  '; ForceText = 0 },
    @{ Addr = 'D11'; Value = '23.18'; ForceText = 1 },
    @{ Addr = 'E11'; Value = '  -0.48%  '; ForceText = 0 },
    @{ Addr = 'D12'; Value = '1.894.91'; ForceText = 1 },
    @{ Addr = 'E12'; Value = '  -1.86%  '; ForceText = 0 },
    @{ Addr = 'D13'; Value = '5.935'; ForceText = 1 },
    @{ Addr = 'E13'; Value = '  -2.16%  '; ForceText = 0 },
    @{ Addr = 'D14'; Value = '7.096'; ForceText = 1 },
    @{ Addr = 'E14'; Value = '  -2.13%  '; ForceText = 0 },
    @{ Addr = 'D15'; Value = '89.04'; ForceText = 1 },
    @{ Addr = 'E15'; Value = '  -2.96%  '; ForceText = 0 },
    @{ Addr = 'D16'; Value = '0.06815'; ForceText = 1 },
    @{ Addr = 'E16'; Value = '  -1.04%  '; ForceText = 0 },
    @{ Addr = 'D17'; Value = '1.005'; ForceText = 1 },
    @{ Addr = 'E17'; Value = '  -0.37%  '; ForceText = 0 },
    @{ Addr = 'D18'; Value = '0.00001019'; ForceText = 1 },
    @{ Addr = 'E18'; Value = '  -1.98%  '; ForceText = 0 },
    @{ Addr = 'D19'; Value = '17.32'; ForceText = 1 },
    @{ Addr = 'E19'; Value = '  -1.71%  '; ForceText = 0 },
    @{ Addr = 'D20'; Value = '1.003'; ForceText = 1 },
    @{ Addr = 'E20'; Value = '  -0.49%  '; ForceText = 0 },
    @{ Addr = 'D21'; Value = '29.506.24'; ForceText = 1 },
    @{ Addr = 'E21'; Value = '  +0.17%  '; ForceText = 0 },
    @{ Addr = 'D22'; Value = '5.506'; ForceText = 1 },
    @{ Addr = 'E22'; Value = '  -2.86%  '; ForceText = 0 },
    @{ Addr = 'D23'; Value = '11.61'; ForceText = 1 },
    @{ Addr = 'E23'; Value = '  -1.10%  '; ForceText = 0 },
    @{ Addr = 'E24'; Value = '  -1.64%  '; ForceText = 0 },
    @{ Addr = 'D25'; Value = '2.144.60'; ForceText = 1 },
    @{ Addr = 'E25'; Value = '  -0.71%  '; ForceText = 0 },
    @{ Addr = 'D26'; Value = '157.07'; ForceText = 1 },
    @{ Addr = 'E26'; Value = '  +0.66%  '; ForceText = 0 },
    @{ Addr = 'D27'; Value = '6.506'; ForceText = 1 },
    @{ Addr = 'E27'; Value = '  -2.20%  '; ForceText = 0 },
    @{ Addr = 'D28'; Value = '19.63'; ForceText = 1 },
    @{ Addr = 'E28'; Value = '  -1.92%  '; ForceText = 0 },
    @{ Addr = 'D29'; Value = '2.055'; ForceText = 1 },
    @{ Addr = 'E29'; Value = '  -2.97%  '; ForceText = 0 },
    @{ Addr = 'D30'; Value = '119.20'; ForceText = 1 },
    @{ Addr = 'E30'; Value = '  -1.43%  '; ForceText = 0 },
    @{ Addr = 'D31'; Value = '0.9943'; ForceText = 1 },
    @{ Addr = 'D32'; Value = '0.09539'; ForceText = 1 },
    @{ Addr = 'E32'; Value = '  -0.77%  '; ForceText = 0 },
    @{ Addr = 'D33'; Value = '5.473'; ForceText = 1 },
    @{ Addr = 'E33'; Value = '  -3.10%  '; ForceText = 0 },
    @{ Addr = 'D34'; Value = '1.388'; ForceText = 1 },
    @{ Addr = 'E34'; Value = '  +0.99%  '; ForceText = 0 },
    @{ Addr = 'D35'; Value = '3.528'; ForceText = 1 },
    @{ Addr = 'E35'; Value = '  -0.46%  '; ForceText = 0 },
    @{ Addr = 'D36'; Value = '0.06382'; ForceText = 1 },
    @{ Addr = 'E36'; Value = '  +4.59%  '; ForceText = 0 },
    @{ Addr = 'D37'; Value = '0.02243'; ForceText = 1 },
    @{ Addr = 'E37'; Value = '  -1.89%  '; ForceText = 0 },
    @{ Addr = 'D38'; Value = '1.197'; ForceText = 1 },
    @{ Addr = 'E38'; Value = '  +1.15%  '; ForceText = 0 },
    @{ Addr = 'D39'; Value = '0.5814'; ForceText = 1 },
    @{ Addr = 'E39'; Value = '  -2.66%  '; ForceText = 0 },
    @{ Addr = 'D40'; Value = '10.54'; ForceText = 1 },
    @{ Addr = 'E40'; Value = '  -3.12%  '; ForceText = 0 },
    @{ Addr = 'D41'; Value = '7.720'; ForceText = 1 },
    @{ Addr = 'E41'; Value = '  -4.38%  '; ForceText = 0 },
    @{ Addr = 'D42'; Value = '0.1817'; ForceText = 1 },
    @{ Addr = 'E42'; Value = '  -1.59%  '; ForceText = 0 },
    @{ Addr = 'D43'; Value = '2.430'; ForceText = 1 },
    @{ Addr = 'E43'; Value = '  +1.62%  '; ForceText = 0 },
    @{ Addr = 'D44'; Value = '1.265'; ForceText = 1 },
    @{ Addr = 'E44'; Value = '  -1.17%  '; ForceText = 0 },
    @{ Addr = 'B45'; Value = 'EnergySwap'; ForceText = 0 },
    @{ Addr = 'C45'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = 0 },
    @{ Addr = 'D45'; Value = '12.18'; ForceText = 1 },
    @{ Addr = 'E45'; Value = '  -1.93%  '; ForceText = 0 },
    @{ Addr = 'B46'; Value = 'Cronos'; ForceText = 0 },
    @{ Addr = 'C46'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; ForceText = 0 },
    @{ Addr = 'D46'; Value = '0.07364'; ForceText = 1 },
    @{ Addr = 'E46'; Value = '  -3.12%  '; ForceText = 0 },
    @{ Addr = 'D47'; Value = '0.5489'; ForceText = 1 },
    @{ Addr = 'E47'; Value = '  -1.89%  '; ForceText = 0 },
    @{ Addr = 'D48'; Value = '1.952'; ForceText = 1 },
    @{ Addr = 'E48'; Value = '  -0.14%  '; ForceText = 0 },
    @{ Addr = 'D49'; Value = '116.20'; ForceText = 1 },
    @{ Addr = 'E49'; Value = '  -1.93%  '; ForceText = 0 },
    @{ Addr = 'D50'; Value = '2.375'; ForceText = 1 },
    @{ Addr = 'E50'; Value = '  -2.14%  '; ForceText = 0 },
    @{ Addr = 'D51'; Value = '71.08'; ForceText = 1 },
    @{ Addr = 'E51'; Value = '  -1.68%  '; ForceText = 0 },
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    if ($u.ForceText -eq 1) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
